$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add two more rows to the "R Assignments" table (rows 20 and 21),
# reusing the same category labels / point values as the Python table
# (Control Structures / Functions, both worth 15 points).
$ws.Range("A20").Value = "Control Structures"
$ws.Range("B20").Value = 15

$ws.Range("A21").Value = "Functions"
$ws.Range("B21").Value = 15

# Keep selection/active cell consistent with the extended range.
$ws.Range("B22").Select()
